$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column before EZ (which currently
#     holds "01-oct."), shifting all subsequent day columns one to the right
#     (EZ..GD -> FA..GE). The new EZ column is the "29-dec" header with "-"
#     placeholder values for every hour row. ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("EZ1").EntireColumn.Insert()

$ws1.Range("EZ1").Value = "29-dec"
$ws1.Range("EZ2:EZ25").Value = "-"

# --- Sheet "Gaz": append two new trailing dates with blank prices. ---
$ws2 = $wb.Worksheets.Item("Gaz")

$ws2.Range("A184").Value = "'2025-12-27"
$ws2.Range("A184").ClearFormats()
$ws2.Range("B184").Value = "'"
$ws2.Range("B184").ClearFormats()

$ws2.Range("A185").Value = "'2025-12-28"
$ws2.Range("A185").ClearFormats()
$ws2.Range("B185").Value = "'"
$ws2.Range("B185").ClearFormats()
